$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 437, pushing the existing rows
# 437..465 down to 438..466 (dimension grows from R465 to R466).
$ws.Rows.Item(437).Insert()

# Populate the newly inserted row 437 with the new record.
$ws.Cells.Item(437, 1).Value2  = 3
$ws.Cells.Item(437, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(437, 3).Value2  = "Coquimbo"
$ws.Cells.Item(437, 4).Value2  = 44826
$ws.Cells.Item(437, 5).Value2  = 5
$ws.Cells.Item(437, 6).Value2  = 100112017
$ws.Cells.Item(437, 7).Value2  = "Apio"
$ws.Cells.Item(437, 8).Value2  = "Americana (o)"
$ws.Cells.Item(437, 9).Value2  = "Primera"
$ws.Cells.Item(437, 10).Value2 = 230
$ws.Cells.Item(437, 11).Value2 = 9000
$ws.Cells.Item(437, 12).Value2 = 10000
$ws.Cells.Item(437, 13).Value2 = 9522
$ws.Cells.Item(437, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(437, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(437, 16).Value2 = 1587
$ws.Cells.Item(437, 17).Value2 = 6
$ws.Cells.Item(437, 18).Value2 = "Hortaliza"
